$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the two runs of the "Meta" paragraph ("Monitoria" + " las
#    condiciones...") into a single run while preserving the first run's
#    (empty) rPr container.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(17)
$metaStart = $metaPara.Range.Start

$run1 = $d.Range($metaStart, $metaStart + 9)   # "Monitoria"

$tailText = " las condiciones de temperatura del " + [char]0x00F3 + "rgano y traslado accidentado. As" + [char]0x00ED + " como evitar el acceso de este por personas no autorizadas. "

$metaEnd = $metaPara.Range.End
$run2 = $d.Range($metaStart + 9, $metaEnd - 1)  # everything up to (not incl.) the paragraph mark
$run2.Delete()

$insPoint = $d.Range($metaStart + 9, $metaStart + 9)
$insPoint.InsertAfter($tailText)

# ---------------------------------------------------------------------------
# 2) Insert the "Investigacion sobre..." block plus the long article text as
#    new paragraphs right after the Meta paragraph. The first of these
#    inherits jc="both" from the Meta paragraph (as in the target), the rest
#    get their alignment explicitly cleared to match plain "Normal" pPr.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(17)
$cur = $metaPara.Range
$cur.Collapse(0)
$cur.InsertParagraphAfter()

$blockA = "Investigaci" + [char]0x00F3 + "n sobre el traslado de " + [char]0x00F3 + "rganos en condiciones ideales."

$restBlock = @(
  "El traslado de " + [char]0x00F3 + "rganos para el trasplante constituye una actividad que se lleva a cabo siempre, bien sea por aire o por carretera, " + [char]0x201C + "in extremis" + [char]0x201D + ", y por lo tanto se trata de una tarea que nunca est" + [char]0x00E1 + " exenta de riesgos.",
  "Precisamente la pasada madrugada fallecieron dos pilotos que viajaban a bordo de una avioneta que parti" + [char]0x00F3 + " de Oporto -donde recogi" + [char]0x00F3 + " un coraz" + [char]0x00F3 + "n para trasplante- y que se estrell" + [char]0x00F3 + " cerca del aeropuerto de Santiago de Compostela tras haber dejado el " + [char]0x00F3 + "rgano en el aeropuerto asturiano situado en las inmediaciones de Avil" + [char]0x00E9 + "s.",
  "Fuentes del Ministerio de Sanidad han explicado a Efe que el protocolo de actuaci" + [char]0x00F3 + "n que se sigue en el caso de " + [char]0x00F3 + "rganos v" + [char]0x00E1 + "lidos para el trasplante es siempre similar.",
  "Procedimiento r" + [char]0x00E1 + "pido y efectivo",
  "El proceso comienza cuando hay un hospital con un posible donante, en cuyo caso se pone en contacto con la Organizaci" + [char]0x00F3 + "n Nacional de Trasplante, que tiene sus l" + [char]0x00ED + "neas abiertas las 24 horas de todos los d" + [char]0x00ED + "as del a" + [char]0x00F1 + "o.",
  "A continuaci" + [char]0x00F3 + "n, la ONT localiza cu" + [char]0x00E1 + "l es el receptor m" + [char]0x00E1 + "s adecuado y, una vez identificado, deciden cu" + [char]0x00E1 + "l es el m" + [char]0x00E9 + "todo de transporte para el " + [char]0x00F3 + "rgano m" + [char]0x00E1 + "s adecuado.",
  "Para distancias cortas se utiliza habitualmente la ambulancia, mientras que el avi" + [char]0x00F3 + "n se emplea para trayectos superiores a las dos horas.",
  "Los aviones pertenecen a empresas que trabajan y colaboran habitualmente con la ONT, y una vez localizado el aparato m" + [char]0x00E1 + "s pr" + [char]0x00F3 + "ximo, el personal de enfermer" + [char]0x00ED + "a agiliza los permisos de aviaci" + [char]0x00F3 + "n civil y otros tr" + [char]0x00E1 + "mites administrativos.",
  "Hasta este momento, el proceso se ha prolongado durante unas dos horas, seg" + [char]0x00FA + "n las mimas fuentes.",
  "El caso de tratarse de un " + [char]0x00F3 + "rgano vital, un equipo m" + [char]0x00E9 + "dico procedente del hospital receptor suele trasladarse al lugar donde se encuentra el donante para participar en la extracci" + [char]0x00F3 + "n y recibir ese " + [char]0x00F3 + "rgano.",
  "Si no se trata de un " + [char]0x00F3 + "rgano vital, ese equipo m" + [char]0x00E9 + "dico espera en el hospital en el que se va a practicar el trasplante.",
  "Investigaci" + [char]0x00F3 + "n del sistema de conversaci" + [char]0x00F3 + "n de temperatura.",
  "Si los " + [char]0x00F3 + "rganos no se conservan apropiadamente, se deterioran en cuesti" + [char]0x00F3 + "n de muy pocas horas. Por eso no puede haber bancos de " + [char]0x00F3 + "rganos y por eso el tiempo es un factor muy valioso.",
  "En funci" + [char]0x00F3 + "n de los " + [char]0x00F3 + "rganos el tiempo de conservaci" + [char]0x00F3 + "n es de las 3 a 5 h del coraz" + [char]0x00F3 + "n o pulm" + [char]0x00F3 + "n,  de 12 a 24 para el h" + [char]0x00ED + "gado y el p" + [char]0x00E1 + "ncreas y 48-72h para los ri" + [char]0x00F1 + "ones.",
  "Para conservar los " + [char]0x00F3 + "rganos s" + [char]0x00F3 + "lidos, se asocian la hipotermia a 4" + [char]0x00BA + "C, y el uso de soluciones con las que se lavan los " + [char]0x00F3 + "rganos, as" + [char]0x00ED + " como con las que se perfunden para que su enfriamiento sea alcanzado de la forma m" + [char]0x00E1 + "s r" + [char]0x00E1 + "pida y homog" + [char]0x00E9 + "nea posible. Las soluciones de preservaci" + [char]0x00F3 + "n pretenden disminuir y frenar todos los procesos de degradaci" + [char]0x00F3 + "n celular y permitir que los " + [char]0x00F3 + "rganos funcionen adecuadamente. Actualmente hay algunas l" + [char]0x00ED + "neas de investigaci" + [char]0x00F3 + "n que han demostrado que pretratando un " + [char]0x00F3 + "rgano con sangre, incluso a temperatura normal o cercana a lo normal antes de implantarlo, en vez de guardarlo exclusivamente en fr" + [char]0x00ED + "o, se podr" + [char]0x00ED + "a mejorar la funci" + [char]0x00F3 + "n postrasplante.",
  "Para el proceso de conservaci" + [char]0x00F3 + "n, se utiliza a veces el sistema circulatorio del paciente, para  perfundir los l" + [char]0x00ED + "quidos y as" + [char]0x00ED + " sacar la sangre del " + [char]0x00F3 + "rgano.",
  "Cuando se extrae el " + [char]0x00F3 + "rgano se introduce en una soluci" + [char]0x00F3 + "n a 8" + [char]0x00BA + "C y luego se revisa. Posteriormente se mete en un recipiente est" + [char]0x00E9 + "ril, que a su vez es introducido en una doble bolsa est" + [char]0x00E9 + "ril, que  a su vez se introduce en un contenedor isotermo que contiene una soluci" + [char]0x00F3 + "n fr" + [char]0x00ED + "a con hielo para conservarlo durante el transporte",
  "Log" + [char]0x00ED + "stica del trasplante",
  "Cuando el coordinador de trasplantes de cualquier centro del Sistema Nacional de Salud, detecta la existencia de un posible donante, se pone en marcha la maquinaria del proceso de donaci" + [char]0x00F3 + "n/trasplante, comunicando dicho donante al coordinador de guardia de la oficina central de la ONT, " + [char]0x00F3 + " a la OCATT si la donaci" + [char]0x00F3 + "n se produce dentro del territorio de Catalu" + [char]0x00F1 + "a. Una vez que se ha diagnosticado la muerte encef" + [char]0x00E1 + "lica, hay que llamar a la ONT aunque no se disponga de los permisos, para ir comunicando a los diferentes equipos las caracter" + [char]0x00ED + "sticas antropom" + [char]0x00E9 + "tricas, y datos anal" + [char]0x00ED + "ticos  del " + [char]0x00F3 + "rgano a trasplantar. Una vez realizado el trasplante el coordinador de trasplante lo comunica a la ONT para as" + [char]0x00ED + " excluirlo de la lista de espera y modificar el turno correspondiente.",
  "Transporte de los equipos de trasplante",
  "1.- Donante local.- Es  aquel que est" + [char]0x00E1 + " en la misma ciudad que el equipo extractor/trasplantador, pero en otro hospital. En estos casos es el coordinador del hospital quien se encarga de organizar el desplazamiento del equipo en funci" + [char]0x00F3 + "n de los acuerdos internos previstos.",
  "2.- Distancias cortas.- Cuando la distancia es inferior a 200 Kms, el traslado de los equipos se realiza preferentemente mediante autom" + [char]0x00F3 + "viles sanitarios " + [char]0x00F3 + " helic" + [char]0x00F3 + "pteros, este " + [char]0x00FA + "ltimo siempre que la climatolog" + [char]0x00ED + "a y el horario lo permita, ya que la mayor" + [char]0x00ED + "a de los helic" + [char]0x00F3 + "pteros no disponen de ayudas nocturnas.",
  "3.- Distancias largas.- En estos casos, y dado el corto tiempo de isquemia fr" + [char]0x00ED + "a (tiempo transcurrido desde la colocaci" + [char]0x00F3 + "n en soluci" + [char]0x00F3 + "n de transporte y el inicio de la desinfecci" + [char]0x00F3 + "n) que toleran los " + [char]0x00F3 + "rganos, se contratan aviones de compa" + [char]0x00F1 + "" + [char]0x00ED + "as privadas de aviaci" + [char]0x00F3 + "n y ocasionalmente se recurre a aviones del Ejercito del Aire. La preparaci" + [char]0x00F3 + "n de un vuelo necesita un tiempo no inferior a 2 hs (verificaci" + [char]0x00F3 + "n del avi" + [char]0x00F3 + "n, aviso a la tripulaci" + [char]0x00F3 + "n, preparaci" + [char]0x00F3 + "n del plan del vuelo, etc.), por eso es importante que se comunique la existencia del donante a la ONT a la mayor brevedad posible. Una vez contratado el vuelo y con los horarios previstos, se avisa a los coordinadores de los hospitales implicados. La mayor" + [char]0x00ED + "a de los aeropuertos nacionales no est" + [char]0x00E1 + "n operativos las 24 h, por lo que el coordinador de la ONT debe tenerlo en cuenta, ya que en caso de no ser de 24 h pondr" + [char]0x00E1 + " en marcha los mecanismos necesarios para su apertura " + [char]0x00F3 + " para que se mantenga operativo fuera de horario."
)

$curParaIndex = $metaPara.Index + 1
$pA = $d.Paragraphs.Item($curParaIndex)
$pA.Range.Text = $blockA

$prevPara = $pA
foreach ($txt in $restBlock) {
    $r = $prevPara.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $curParaIndex = $prevPara.Index + 1
    $newPara = $d.Paragraphs.Item($curParaIndex)
    $newPara.Range.ParagraphFormat.Alignment = 0
    $newPara.Range.Text = $txt
    $prevPara = $newPara
}

# ---------------------------------------------------------------------------
# 3) The paragraph that used to be the 3rd trailing empty paragraph (the one
#    carrying widowControl/bidi/spacing/jc=left) now gets the
#    "Investigacion de mecanismo..." text typed directly into its run.
# ---------------------------------------------------------------------------
$lastEmptyIndex = $prevPara.Index + 3
$doorPara = $d.Paragraphs.Item($lastEmptyIndex)
$doorPara.Range.Text = "Investigaci" + [char]0x00F3 + "n de mecanismo de apertura y cierre de puerta."

# ---------------------------------------------------------------------------
# 4) Append the remaining "puertas automaticas" article paragraphs after it;
#    they all inherit the widowControl/bidi/spacing/jc=left pPr.
# ---------------------------------------------------------------------------
$doorBlock = @(
  "Investigaci" + [char]0x00F3 + "n del sistema de monitoreo y control de acceso.",
  "La evoluci" + [char]0x00F3 + "n de la tecnolog" + [char]0x00ED + "a no ha dejado rinc" + [char]0x00F3 + "n sin influenciar, es cierto que lo m" + [char]0x00E1 + "s probable, cuando mencionamos " + [char]0x201C + "aparatos modernos" + [char]0x201D + ", lo primero que se nos cruce por la cabeza sea una notebook, un reproductor de mp3 o un televisor digital. Pero existen otros dispositivos que son mucho m" + [char]0x00E1 + "s comunes y que tambi" + [char]0x00E9 + "n han sido innovados por los avances tecnol" + [char]0x00F3 + "gicos, nos referimos a las puertas autom" + [char]0x00E1 + "ticas; una simple puerta era un objeto que oper" + [char]0x00E1 + "bamos manualmente cada vez que quer" + [char]0x00ED + "amos salir o entrar en un cuarto o dependencia, hoy esta simple acci" + [char]0x00F3 + "n se ha vuelto incluso m" + [char]0x00E1 + "s sencilla mediante un sistema de automatizaci" + [char]0x00F3 + "n.",
  "Podemos se" + [char]0x00F1 + "alar que las primeras puertas autom" + [char]0x00E1 + "ticas fueron aquellas destinadas a los asensores, justamente lo que se intent" + [char]0x00F3 + " con el dise" + [char]0x00F1 + "o de estas puertas era hacer el trabajo de esta m" + [char]0x00E1 + "quina mucho m" + [char]0x00E1 + "s simple evitando que las mismas permanezcan sin funcionar debido a que una persona cerr" + [char]0x00F3 + " incorrectamente alguna de las puertas en cuesti" + [char]0x00F3 + "n.",
  "as puertas autom" + [char]0x00E1 + "ticas modernas sustituyen a las antiguas de tipo manual pero la ventaja es que para remplazar estas " + [char]0x00FA + "ltimas no se necesita de obra manual; las puertas manuales de los asensores dificultan el acceso a personas que tienen sus manos ocupadas o sufren de alguna discapacidad motriz y es por eso que precisan una apertura autom" + [char]0x00E1 + "tica.",
  "Este problema queda resuelto mediante la instalaci" + [char]0x00F3 + "n de una puerta autom" + [char]0x00E1 + "tica que se adapta a las necesidades del individuo que usa el elevador; actualmente los dise" + [char]0x00F1 + "os de estos elevadores son telesc" + [char]0x00F3 + "picos de 4 hojas de apertura central; las puertas se abren y cierran mediante u operador de alto rendimiento y funcionamiento silencioso, las ventajas m" + [char]0x00E1 + "s significativas radican en la fiabilidad, funcionamiento, seguridad y est" + [char]0x00E9 + "tica."
)

$prevPara = $doorPara
foreach ($txt in $doorBlock) {
    $r = $prevPara.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $curParaIndex = $prevPara.Index + 1
    $newPara = $d.Paragraphs.Item($curParaIndex)
    $newPara.Range.Text = $txt
    $prevPara = $newPara
}

# ---------------------------------------------------------------------------
# 5) Append one final empty paragraph (still widowControl style) whose run
#    and pPr carry single-underline character formatting.
# ---------------------------------------------------------------------------
$r = $prevPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$curParaIndex = $prevPara.Index + 1
$finalPara = $d.Paragraphs.Item($curParaIndex)
$finalPara.Range.Font.Underline = 1

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
